$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.119.41"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "1.835.93"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.18"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6643"
$ws.Range("E6").Value = "  -4.53%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2953"
$ws.Range("E8").Value = "  -4.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07357"
$ws.Range("E9").Value = "  -4.43%  "

$ws.Range("E10").Value = "  -3.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07679"
$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "1.840.19"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.024"
$ws.Range("E13").Value = "  -2.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6761"
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.22"
$ws.Range("E15").Value = "  -5.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.220"
$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("D17").Value = "29.045.40"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008234"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "229.03"
$ws.Range("E19").Value = "  -3.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.52"
$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.311"

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.89"
$ws.Range("E24").Value = "  +0.53%  "

$ws.Range("E25").Value = "  -5.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.690"
$ws.Range("E26").Value = "  -2.44%  "

$ws.Range("E27").Value = "  -1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.502"
$ws.Range("E28").Value = "  -2.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.228"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.103"
$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("E31").Value = "  -1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05300"
$ws.Range("E32").Value = "  +3.69%  "

$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7497"
$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.130"
$ws.Range("E35").Value = "  -1.64%  "

$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").Value = "1.313.40"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  -3.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.716"
$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9223"
$ws.Range("E40").Value = "  -2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.985"
$ws.Range("E41").Value = "  +3.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9984"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.69"
$ws.Range("E43").Value = "  -1.98%  "

$ws.Range("D44").Value = "1.986.68"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5162"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("E46").Value = "  -2.91%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.763"
$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.70"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.284"
$ws.Range("E49").Value = "  -5.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05936"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07245"
$ws.Range("E51").Value = "  +7.53%  "
